$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell X1: "NaN" placeholder -> new region header "ColRisaralda" ---
# (adds a new shared string; existing string table reused for all "NaN" writes below)
$ws.Cells.Item(1, 24).Value = "ColRisaralda"

# --- Scattered historical-data corrections: some cells flip between a numeric
#     value and the literal text "NaN", others are simple numeric decrements ---
$ws.Cells.Item(16, 18).Value = 2  # R16: was NaN -> 2
$ws.Cells.Item(18, 69).Value = 1  # BQ18: was NaN -> 1
$ws.Cells.Item(19, 29).Value = "NaN"  # AC19: was 1 -> NaN
$ws.Cells.Item(23, 89).Value = "NaN"  # CK23: was 1 -> NaN
$ws.Cells.Item(25, 69).Value = "NaN"  # BQ25: was 4 -> NaN
$ws.Cells.Item(26, 73).Value = "NaN"  # BU26: was 20 -> NaN
$ws.Cells.Item(27, 101).Value = "NaN"  # CW27: was 1 -> NaN
$ws.Cells.Item(28, 29).Value = 2  # AC28: was NaN -> 2
$ws.Cells.Item(30, 91).Value = "NaN"  # CM30: was 1 -> NaN
$ws.Cells.Item(31, 91).Value = "NaN"  # CM31: was 1 -> NaN
$ws.Cells.Item(31, 106).Value = 1  # DB31: was NaN -> 1
$ws.Cells.Item(36, 37).Value = "NaN"  # AK36: was 2 -> NaN
$ws.Cells.Item(37, 37).Value = "NaN"  # AK37: was 2 -> NaN
$ws.Cells.Item(37, 73).Value = "NaN"  # BU37: was 33 -> NaN
$ws.Cells.Item(41, 29).Value = "NaN"  # AC41: was 9 -> NaN
$ws.Cells.Item(72, 24).Value = "NaN"  # X72: was 3 -> NaN
$ws.Cells.Item(73, 24).Value = "NaN"  # X73: was 3 -> NaN
$ws.Cells.Item(83, 91).Value = 25  # CM83: was NaN -> 25
$ws.Cells.Item(88, 121).Value = "NaN"  # DQ88: was 1 -> NaN
$ws.Cells.Item(91, 112).Value = "NaN"  # DH91: was 2 -> NaN
$ws.Cells.Item(92, 112).Value = "NaN"  # DH92: was 2 -> NaN
$ws.Cells.Item(105, 73).Value = "NaN"  # BU105: was 105 -> NaN
$ws.Cells.Item(106, 29).Value = 41  # AC106: was NaN -> 41
$ws.Cells.Item(107, 73).Value = "NaN"  # BU107: was 108 -> NaN
$ws.Cells.Item(108, 69).Value = "NaN"  # BQ108: was 467 -> NaN
$ws.Cells.Item(126, 85).Value = "NaN"  # CG126: was 72 -> NaN
$ws.Cells.Item(127, 90).Value = "NaN"  # CL127: was 187 -> NaN
$ws.Cells.Item(141, 90).Value = 350  # CL141: was 351 -> 350
$ws.Cells.Item(142, 90).Value = 370  # CL142: was 371 -> 370
$ws.Cells.Item(143, 90).Value = 373  # CL143: was 374 -> 373
$ws.Cells.Item(144, 90).Value = 386  # CL144: was 388 -> 386
$ws.Cells.Item(145, 90).Value = 393  # CL145: was 396 -> 393
$ws.Cells.Item(146, 90).Value = 416  # CL146: was 420 -> 416
$ws.Cells.Item(147, 90).Value = 440  # CL147: was 445 -> 440
$ws.Cells.Item(148, 90).Value = 497  # CL148: was 502 -> 497
$ws.Cells.Item(149, 90).Value = 508  # CL149: was 510 -> 508
$ws.Cells.Item(150, 90).Value = 540  # CL150: was 542 -> 540
$ws.Cells.Item(150, 91).Value = "NaN"  # CM150: was 464 -> NaN
$ws.Cells.Item(151, 90).Value = "NaN"  # CL151: was 578 -> NaN
$ws.Cells.Item(151, 91).Value = "NaN"  # CM151: was 486 -> NaN
$ws.Cells.Item(152, 90).Value = 591  # CL152: was 592 -> 591
$ws.Cells.Item(153, 90).Value = 603  # CL153: was 604 -> 603
$ws.Cells.Item(153, 99).Value = "NaN"  # CU153: was 401 -> NaN
$ws.Cells.Item(154, 90).Value = 621  # CL154: was 622 -> 621
$ws.Cells.Item(154, 99).Value = "NaN"  # CU154: was 417 -> NaN
$ws.Cells.Item(155, 32).Value = "NaN"  # AF155: was 176 -> NaN
$ws.Cells.Item(156, 90).Value = 668  # CL156: was 669 -> 668
$ws.Cells.Item(157, 69).Value = 3048  # BQ157: was 3050 -> 3048
$ws.Cells.Item(157, 90).Value = 685  # CL157: was 687 -> 685
$ws.Cells.Item(158, 69).Value = 3119  # BQ158: was 3121 -> 3119
$ws.Cells.Item(158, 90).Value = 697  # CL158: was 699 -> 697
$ws.Cells.Item(159, 69).Value = 3156  # BQ159: was 3158 -> 3156
$ws.Cells.Item(159, 90).Value = 715  # CL159: was 717 -> 715
$ws.Cells.Item(160, 69).Value = 3404  # BQ160: was 3406 -> 3404
$ws.Cells.Item(162, 90).Value = 827  # CL162: was 828 -> 827
$ws.Cells.Item(163, 69).Value = 3765  # BQ163: was 3767 -> 3765
$ws.Cells.Item(164, 69).Value = 3882  # BQ164: was 3884 -> 3882
$ws.Cells.Item(165, 69).Value = 4108  # BQ165: was 4110 -> 4108
$ws.Cells.Item(165, 90).Value = 948  # CL165: was 949 -> 948
$ws.Cells.Item(166, 69).Value = 4246  # BQ166: was 4248 -> 4246
$ws.Cells.Item(166, 90).Value = 963  # CL166: was 964 -> 963

# --- Append newly-reported day: row 169 (2020-08-20 series) ---
$vals169 = @(44063,513719,2677,66613,62271,179540,22615,2416,2024,4428,3703,6968,3571,16351,17924,4006,2901,10815,5287,12091,8103,2282,716,3865,11806,9827,5142,41259,780,115,165,436,28,15,205,1920,2022,35127,5527,2351,32537,775,18631,1390,5672,1352,1527,2923,1318,925,2445,2559,38124,10546,1621,6342,2480,273,1360,2491,724,1884,7346,7075,6846,13438,1838,732,5023,4407,4849,1074,1246,1900,2216,535,3621,2056,864,583,1608,1641,788,696,3724,928,995,1005,1320,1183,1114,956,927,1011,474,2759,781,718,629,1059,954,525,646,710,992,814,900,720,309,319,611,486,365,513,299,510,689,497,466,343,508,106675,214931,6878,92781,60624,19958,6738)
for ($i = 0; $i -lt $vals169.Length; $i++) {
    $ws.Cells.Item(169, $i + 1).Value = $vals169[$i]
}

# --- Restore frozen-pane selection to the top of the scrollable area ---
$ws.Range("B2").Select()
